$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values as text even when they
# look like numbers (Excel auto-converts numeric-looking strings when
# assigned through .Value). Temporarily force text format on the whole
# data range, write the values, then restore the default "Normal" style
# so the saved file does not carry a permanent number-format change.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.517.49"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").Value = "1.680.58"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "216.67"
$ws.Range("E5").Value = "  +2.75%  "

$ws.Range("D6").Value = "0.5325"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.2681"
$ws.Range("E8").Value = "  +3.94%  "

$ws.Range("D9").Value = "0.06398"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("D10").Value = "21.66"
$ws.Range("E10").Value = "  +5.69%  "

$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  +3.01%  "

$ws.Range("D12").Value = "1.675.60"
$ws.Range("E12").Value = "  +2.75%  "

$ws.Range("D13").Value = "4.502"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("D14").Value = "0.5571"
$ws.Range("E14").Value = "  +1.26%  "

$ws.Range("D15").Value = "0.0₅8338"
$ws.Range("E15").Value = "  +4.12%  "

$ws.Range("D16").Value = "65.72"
$ws.Range("E16").Value = "  +1.57%  "

$ws.Range("D17").Value = "26.558.79"
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "4.765"
$ws.Range("E19").Value = "  +1.97%  "

$ws.Range("D20").Value = "194.68"
$ws.Range("E20").Value = "  +4.93%  "

$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("D22").Value = "6.344"
$ws.Range("E22").Value = "  +3.90%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "143.25"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").Value = "0.1284"
$ws.Range("E25").Value = "  +5.79%  "

$ws.Range("D26").Value = "7.435"
$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("D27").Value = "16.35"
$ws.Range("E27").Value = "  +4.24%  "

$ws.Range("D28").Value = "1.428"
$ws.Range("E28").Value = "  +4.79%  "

$ws.Range("D29").Value = "0.06168"
$ws.Range("E29").Value = "  +4.63%  "

$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").Value = "3.606"
$ws.Range("E31").Value = "  +5.43%  "

$ws.Range("D32").Value = "3.456"
$ws.Range("E32").Value = "  +2.00%  "

$ws.Range("D33").Value = "1.689"
$ws.Range("E33").Value = "  +4.03%  "

$ws.Range("D34").Value = "1.008"
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").Value = "2.423"
$ws.Range("E35").Value = "  +1.73%  "

$ws.Range("D36").Value = "2.790"
$ws.Range("E36").Value = "  +2.37%  "

$ws.Range("D37").Value = "0.5730"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").Value = "0.01638"
$ws.Range("E38").Value = "  +2.26%  "

$ws.Range("D39").Value = "6.025"
$ws.Range("E39").Value = "  +6.18%  "

$ws.Range("D40").Value = "1.075.32"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("D41").Value = "0.8592"
$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").Value = "100.05"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "1.826.79"
$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("D45").Value = "0.0₈110"
$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("D46").Value = "57.01"
$ws.Range("E46").Value = "  +3.61%  "

$ws.Range("D47").Value = "8.118"
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").Value = "0.05210"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").Value = "1.469"
$ws.Range("E50").Value = "  +5.97%  "

$ws.Range("D51").Value = "6.023"
$ws.Range("E51").Value = "  +3.05%  "

$priceRange.Style = "Normal"
